$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.408429333333333
$ws.Range("H2").Value = 13.225288
$ws.Range("I2").Value = 0.05980478019486075
$ws.Range("J2").Value = 0.05980478019486075
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 36.71344366666667
$ws.Range("N2").Value = 110.140331
$ws.Range("O2").Value = 0.2081992981130139
$ws.Range("P2").Value = 0.2081992981130138
$ws.Range("Q2").Value = 161.8486219878142
$ws.Range("R2").Value = 1456.637597890328
$ws.Range("S2").Value = 0.01245131326037308
$ws.Range("T2").Value = 0.01245131326037308

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.408429333333333
$ws.Range("H3").Value = 13.225288
$ws.Range("I3").Value = 0.05980478019486075
$ws.Range("J3").Value = 0.05980478019486075
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.87708033333332
$ws.Range("N3").Value = 266.631241
$ws.Range("O3").Value = 0.5040155293450301
$ws.Range("P3").Value = 0.50401552934503
$ws.Range("Q3").Value = 391.8083280024898
$ws.Range("R3").Value = 3526.274952022407
$ws.Range("S3").Value = 0.03014253794727592
$ws.Range("T3").Value = 0.03014253794727591

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.408429333333333
$ws.Range("H4").Value = 13.225288
$ws.Range("I4").Value = 0.05980478019486075
$ws.Range("J4").Value = 0.05980478019486075
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.22727166666667
$ws.Range("N4").Value = 42.681815
$ws.Range("O4").Value = 0.08068183420648613
$ws.Range("P4").Value = 0.08068183420648613
$ws.Range("Q4").Value = 62.71992174863556
$ws.Range("R4").Value = 564.47929573772
$ws.Range("S4").Value = 0.0048251593604371
$ws.Range("T4").Value = 0.0048251593604371

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.408429333333333
$ws.Range("H5").Value = 13.225288
$ws.Range("I5").Value = 0.05980478019486075
$ws.Range("J5").Value = 0.05980478019486075
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 36.52018433333333
$ws.Range("N5").Value = 109.560553
$ws.Range("O5").Value = 0.20710333833547
$ws.Range("P5").Value = 0.2071033383354699
$ws.Range("Q5").Value = 160.9966518738071
$ws.Range("R5").Value = 1448.969866864264
$ws.Range("S5").Value = 0.01238576962677466
$ws.Range("T5").Value = 0.01238576962677466

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 56.15338233333333
$ws.Range("H6").Value = 168.460147
$ws.Range("I6").Value = 0.7617771395926449
$ws.Range("J6").Value = 0.7617771395926448
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 36.71344366666667
$ws.Range("N6").Value = 110.140331
$ws.Range("O6").Value = 0.2081992981130139
$ws.Range("P6").Value = 0.2081992981130138
$ws.Range("Q6").Value = 2061.584038987629
$ws.Range("R6").Value = 18554.25635088866
$ws.Range("S6").Value = 0.158601465781728
$ws.Range("T6").Value = 0.158601465781728

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 56.15338233333333
$ws.Range("H7").Value = 168.460147
$ws.Range("I7").Value = 0.7617771395926449
$ws.Range("J7").Value = 0.7617771395926448
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 88.87708033333332
$ws.Range("N7").Value = 266.631241
$ws.Range("O7").Value = 0.5040155293450301
$ws.Range("P7").Value = 0.50401552934503
$ws.Range("Q7").Value = 4990.748672628047
$ws.Range("R7").Value = 44916.73805365243
$ws.Range("S7").Value = 0.3839475082547298
$ws.Range("T7").Value = 0.3839475082547297

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 56.15338233333333
$ws.Range("H8").Value = 168.460147
$ws.Range("I8").Value = 0.7617771395926449
$ws.Range("J8").Value = 0.7617771395926448
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.22727166666667
$ws.Range("N8").Value = 42.681815
$ws.Range("O8").Value = 0.08068183420648613
$ws.Range("P8").Value = 0.08068183420648613
$ws.Range("Q8").Value = 798.9094254585339
$ws.Range("R8").Value = 7190.184829126805
$ws.Range("S8").Value = 0.06146157687890502
$ws.Range("T8").Value = 0.06146157687890501

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 56.15338233333333
$ws.Range("H9").Value = 168.460147
$ws.Range("I9").Value = 0.7617771395926449
$ws.Range("J9").Value = 0.7617771395926448
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 36.52018433333333
$ws.Range("N9").Value = 109.560553
$ws.Range("O9").Value = 0.20710333833547
$ws.Range("P9").Value = 0.2071033383354699
$ws.Range("Q9").Value = 2050.731873753477
$ws.Range("R9").Value = 18456.58686378129
$ws.Range("S9").Value = 0.1577665886772821
$ws.Range("T9").Value = 0.157766588677282

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.790021666666667
$ws.Range("H10").Value = 8.370065
$ws.Range("I10").Value = 0.03784945156141002
$ws.Range("J10").Value = 0.03784945156141001
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 36.71344366666667
$ws.Range("N10").Value = 110.140331
$ws.Range("O10").Value = 0.2081992981130139
$ws.Range("P10").Value = 0.2081992981130138
$ws.Range("Q10").Value = 102.4313032879461
$ws.Range("R10").Value = 921.8817295915151
$ws.Range("S10").Value = 0.007880229249048083
$ws.Range("T10").Value = 0.00788022924904808

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.790021666666667
$ws.Range("H11").Value = 8.370065
$ws.Range("I11").Value = 0.03784945156141002
$ws.Range("J11").Value = 0.03784945156141001
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 88.87708033333332
$ws.Range("N11").Value = 266.631241
$ws.Range("O11").Value = 0.5040155293450301
$ws.Range("P11").Value = 0.50401552934503
$ws.Range("Q11").Value = 247.9689798000739
$ws.Range("R11").Value = 2231.720818200665
$ws.Range("S11").Value = 0.01907671136414315
$ws.Range("T11").Value = 0.01907671136414314

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.790021666666667
$ws.Range("H12").Value = 8.370065
$ws.Range("I12").Value = 0.03784945156141002
$ws.Range("J12").Value = 0.03784945156141001
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 14.22727166666667
$ws.Range("N12").Value = 42.681815
$ws.Range("O12").Value = 0.08068183420648613
$ws.Range("P12").Value = 0.08068183420648613
$ws.Range("Q12").Value = 39.69439620755278
$ws.Range("R12").Value = 357.249565867975
$ws.Range("S12").Value = 0.003053763175684111
$ws.Range("T12").Value = 0.00305376317568411

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.790021666666667
$ws.Range("H13").Value = 8.370065
$ws.Range("I13").Value = 0.03784945156141002
$ws.Range("J13").Value = 0.03784945156141001
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 36.52018433333333
$ws.Range("N13").Value = 109.560553
$ws.Range("O13").Value = 0.20710333833547
$ws.Range("P13").Value = 0.2071033383354699
$ws.Range("Q13").Value = 101.8921055606606
$ws.Range("R13").Value = 917.0289500459451
$ws.Range("S13").Value = 0.007838747772534682
$ws.Range("T13").Value = 0.007838747772534679

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.36182833333333
$ws.Range("H14").Value = 31.085485
$ws.Range("I14").Value = 0.1405686286510843
$ws.Range("J14").Value = 0.1405686286510843
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 36.71344366666667
$ws.Range("N14").Value = 110.140331
$ws.Range("O14").Value = 0.2081992981130139
$ws.Range("P14").Value = 0.2081992981130138
$ws.Range("Q14").Value = 380.4184007995039
$ws.Range("R14").Value = 3423.765607195535
$ws.Range("S14").Value = 0.02926628982186464
$ws.Range("T14").Value = 0.02926628982186463

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.36182833333333
$ws.Range("H15").Value = 31.085485
$ws.Range("I15").Value = 0.1405686286510843
$ws.Range("J15").Value = 0.1405686286510843
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 88.87708033333332
$ws.Range("N15").Value = 266.631241
$ws.Range("O15").Value = 0.5040155293450301
$ws.Range("P15").Value = 0.50401552934503
$ws.Range("Q15").Value = 920.929049181876
$ws.Range("R15").Value = 8288.361442636884
$ws.Range("S15").Value = 0.07084877177888123
$ws.Range("T15").Value = 0.07084877177888119

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.36182833333333
$ws.Range("H16").Value = 31.085485
$ws.Range("I16").Value = 0.1405686286510843
$ws.Range("J16").Value = 0.1405686286510843
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 14.22727166666667
$ws.Range("N16").Value = 42.681815
$ws.Range("O16").Value = 0.08068183420648613
$ws.Range("P16").Value = 0.08068183420648613
$ws.Range("Q16").Value = 147.4205466616972
$ws.Range("R16").Value = 1326.784919955275
$ws.Range("S16").Value = 0.0113413347914599
$ws.Range("T16").Value = 0.0113413347914599

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.36182833333333
$ws.Range("H17").Value = 31.085485
$ws.Range("I17").Value = 0.1405686286510843
$ws.Range("J17").Value = 0.1405686286510843
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 36.52018433333333
$ws.Range("N17").Value = 109.560553
$ws.Range("O17").Value = 0.20710333833547
$ws.Range("P17").Value = 0.2071033383354699
$ws.Range("Q17").Value = 378.4158807636894
$ws.Range("R17").Value = 3405.742926873205
$ws.Range("S17").Value = 0.02911223225887855
$ws.Range("T17").Value = 0.02911223225887854
